$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 89.35833500000001
$ws.Range("H2").Value = 268.075005
$ws.Range("I2").Value = 0.9624640326757887
$ws.Range("J2").Value = 0.9624640326757889
$ws.Range("M2").Value = 1.780241
$ws.Range("N2").Value = 5.340723000000001
$ws.Range("O2").Value = 0.2571704734300857
$ws.Range("P2").Value = 0.2571704734300857
$ws.Range("Q2").Value = 159.079371658735
$ws.Range("R2").Value = 1431.714344928615
$ws.Range("S2").Value = 0.2475173309426621
$ws.Range("T2").Value = 0.2475173309426621
$ws.Range("G3").Value = 89.35833500000001
$ws.Range("H3").Value = 268.075005
$ws.Range("I3").Value = 0.9624640326757887
$ws.Range("J3").Value = 0.9624640326757889
$ws.Range("O3").Value = 0.4607709215973151
$ws.Range("P3").Value = 0.4607709215973152
$ws.Range("Q3").Value = 285.02163451608
$ws.Range("R3").Value = 2565.194710644721
$ws.Range("S3").Value = 0.4434754393402916
$ws.Range("T3").Value = 0.4434754393402918
$ws.Range("G4").Value = 89.35833500000001
$ws.Range("H4").Value = 268.075005
$ws.Range("I4").Value = 0.9624640326757887
$ws.Range("J4").Value = 0.9624640326757889
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.4746316666666666
$ws.Range("N4").Value = 1.423895
$ws.Range("O4").Value = 0.06856445302718973
$ws.Range("P4").Value = 0.06856445302718973
$ws.Range("Q4").Value = 42.41229547160833
$ws.Range("R4").Value = 381.710659244475
$ws.Range("S4").Value = 0.06599081995875872
$ws.Range("T4").Value = 0.06599081995875873
$ws.Range("G5").Value = 89.35833500000001
$ws.Range("H5").Value = 268.075005
$ws.Range("I5").Value = 0.9624640326757887
$ws.Range("J5").Value = 0.9624640326757889
$ws.Range("M5").Value = 0.9303213333333334
$ws.Range("N5").Value = 2.790964
$ws.Range("O5").Value = 0.1343925781596098
$ws.Range("P5").Value = 0.1343925781596098
$ws.Range("Q5").Value = 83.13196536164669
$ws.Range("R5").Value = 748.1876882548202
$ws.Range("S5").Value = 0.1293480227371942
$ws.Range("T5").Value = 0.1293480227371942
$ws.Range("G6").Value = 89.35833500000001
$ws.Range("H6").Value = 268.075005
$ws.Range("I6").Value = 0.9624640326757887
$ws.Range("J6").Value = 0.9624640326757889
$ws.Range("M6").Value = 0.547574
$ws.Range("N6").Value = 1.642722
$ws.Range("O6").Value = 0.07910157378579964
$ws.Range("P6").Value = 0.07910157378579964
$ws.Range("Q6").Value = 48.93030092929001
$ws.Range("R6").Value = 440.3727083636101
$ws.Range("S6").Value = 0.07613241969688218
$ws.Range("T6").Value = 0.0761324196968822
$ws.Range("I7").Value = 0.001854741667334279
$ws.Range("J7").Value = 0.001854741667334279
$ws.Range("M7").Value = 1.780241
$ws.Range("N7").Value = 5.340723000000001
$ws.Range("O7").Value = 0.2571704734300857
$ws.Range("P7").Value = 0.2571704734300857
$ws.Range("Q7").Value = 0.3065580936136668
$ws.Range("R7").Value = 2.759022842523001
$ws.Range("S7").Value = 0.000476984792678863
$ws.Range("T7").Value = 0.0004769847926788631
$ws.Range("I8").Value = 0.001854741667334279
$ws.Range("J8").Value = 0.001854741667334279
$ws.Range("O8").Value = 0.4607709215973151
$ws.Range("P8").Value = 0.4607709215973152
$ws.Range("S8").Value = 0.0008546110273825566
$ws.Range("T8").Value = 0.0008546110273825569
$ws.Range("I9").Value = 0.001854741667334279
$ws.Range("J9").Value = 0.001854741667334279
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.4746316666666666
$ws.Range("N9").Value = 1.423895
$ws.Range("O9").Value = 0.06856445302718973
$ws.Range("P9").Value = 0.06856445302718973
$ws.Range("Q9").Value = 0.08173173121055556
$ws.Range("R9").Value = 0.735585580895
$ws.Range("S9").Value = 0.0001271693479275127
$ws.Range("T9").Value = 0.0001271693479275128
$ws.Range("I10").Value = 0.001854741667334279
$ws.Range("J10").Value = 0.001854741667334279
$ws.Range("M10").Value = 0.9303213333333334
$ws.Range("N10").Value = 2.790964
$ws.Range("O10").Value = 0.1343925781596098
$ws.Range("P10").Value = 0.1343925781596098
$ws.Range("Q10").Value = 0.1602016437071112
$ws.Range("R10").Value = 1.441814793364
$ws.Range("S10").Value = 0.0002492635144931071
$ws.Range("T10").Value = 0.0002492635144931071
$ws.Range("I11").Value = 0.001854741667334279
$ws.Range("J11").Value = 0.001854741667334279
$ws.Range("M11").Value = 0.547574
$ws.Range("N11").Value = 1.642722
$ws.Range("O11").Value = 0.07910157378579964
$ws.Range("P11").Value = 0.07910157378579964
$ws.Range("Q11").Value = 0.0942924253246667
$ws.Range("R11").Value = 0.8486318279220002
$ws.Range("S11").Value = 0.0001467129848522395
$ws.Range("T11").Value = 0.0001467129848522395
$ws.Range("G12").Value = 1.963978
$ws.Range("H12").Value = 5.891934
$ws.Range("I12").Value = 0.02115368628977398
$ws.Range("J12").Value = 0.02115368628977398
$ws.Range("M12").Value = 1.780241
$ws.Range("N12").Value = 5.340723000000001
$ws.Range("O12").Value = 0.2571704734300857
$ws.Range("P12").Value = 0.2571704734300857
$ws.Range("Q12").Value = 3.496354158698
$ws.Range("R12").Value = 31.467187428282
$ws.Range("S12").Value = 0.005440103517932687
$ws.Range("T12").Value = 0.005440103517932688
$ws.Range("G13").Value = 1.963978
$ws.Range("H13").Value = 5.891934
$ws.Range("I13").Value = 0.02115368628977398
$ws.Range("J13").Value = 0.02115368628977398
$ws.Range("O13").Value = 0.4607709215973151
$ws.Range("P13").Value = 0.4607709215973152
$ws.Range("Q13").Value = 6.264398499744
$ws.Range("R13").Value = 56.379586497696
$ws.Range("S13").Value = 0.009747003526919646
$ws.Range("T13").Value = 0.009747003526919647
$ws.Range("G14").Value = 1.963978
$ws.Range("H14").Value = 5.891934
$ws.Range("I14").Value = 0.02115368628977398
$ws.Range("J14").Value = 0.02115368628977398
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.4746316666666666
$ws.Range("N14").Value = 1.423895
$ws.Range("O14").Value = 0.06856445302718973
$ws.Range("P14").Value = 0.06856445302718973
$ws.Range("Q14").Value = 0.9321661514366666
$ws.Range("R14").Value = 8.389495362929999
$ws.Range("S14").Value = 0.001450390929967115
$ws.Range("T14").Value = 0.001450390929967116
$ws.Range("G15").Value = 1.963978
$ws.Range("H15").Value = 5.891934
$ws.Range("I15").Value = 0.02115368628977398
$ws.Range("J15").Value = 0.02115368628977398
$ws.Range("M15").Value = 0.9303213333333334
$ws.Range("N15").Value = 2.790964
$ws.Range("O15").Value = 0.1343925781596098
$ws.Range("P15").Value = 0.1343925781596098
$ws.Range("Q15").Value = 1.827130631597333
$ws.Range("R15").Value = 16.444175684376
$ws.Range("S15").Value = 0.002842898438062316
$ws.Range("T15").Value = 0.002842898438062316
$ws.Range("G16").Value = 1.963978
$ws.Range("H16").Value = 5.891934
$ws.Range("I16").Value = 0.02115368628977398
$ws.Range("J16").Value = 0.02115368628977398
$ws.Range("M16").Value = 0.547574
$ws.Range("N16").Value = 1.642722
$ws.Range("O16").Value = 0.07910157378579964
$ws.Range("P16").Value = 0.07910157378579964
$ws.Range("Q16").Value = 1.075423289372
$ws.Range("R16").Value = 9.678809604348
$ws.Range("S16").Value = 0.001673289876892215
$ws.Range("T16").Value = 0.001673289876892215
$ws.Range("G17").Value = 0.16825
$ws.Range("H17").Value = 0.50475
$ws.Range("I17").Value = 0.001812193272151965
$ws.Range("J17").Value = 0.001812193272151965
$ws.Range("M17").Value = 1.780241
$ws.Range("N17").Value = 5.340723000000001
$ws.Range("O17").Value = 0.2571704734300857
$ws.Range("P17").Value = 0.2571704734300857
$ws.Range("Q17").Value = 0.2995255482500001
$ws.Range("R17").Value = 2.695729934250001
$ws.Range("S17").Value = 0.000466042601746137
$ws.Range("T17").Value = 0.000466042601746137
$ws.Range("G18").Value = 0.16825
$ws.Range("H18").Value = 0.50475
$ws.Range("I18").Value = 0.001812193272151965
$ws.Range("J18").Value = 0.001812193272151965
$ws.Range("O18").Value = 0.4607709215973151
$ws.Range("P18").Value = 0.4607709215973152
$ws.Range("Q18").Value = 0.536658276
$ws.Range("R18").Value = 4.829924484
$ws.Range("S18").Value = 0.000835005964121915
$ws.Range("T18").Value = 0.0008350059641219152
$ws.Range("G19").Value = 0.16825
$ws.Range("H19").Value = 0.50475
$ws.Range("I19").Value = 0.001812193272151965
$ws.Range("J19").Value = 0.001812193272151965
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.4746316666666666
$ws.Range("N19").Value = 1.423895
$ws.Range("O19").Value = 0.06856445302718973
$ws.Range("P19").Value = 0.06856445302718973
$ws.Range("Q19").Value = 0.07985677791666666
$ws.Range("R19").Value = 0.71871100125
$ws.Range("S19").Value = 0.0001242520404846527
$ws.Range("T19").Value = 0.0001242520404846527
$ws.Range("G20").Value = 0.16825
$ws.Range("H20").Value = 0.50475
$ws.Range("I20").Value = 0.001812193272151965
$ws.Range("J20").Value = 0.001812193272151965
$ws.Range("M20").Value = 0.9303213333333334
$ws.Range("N20").Value = 2.790964
$ws.Range("O20").Value = 0.1343925781596098
$ws.Range("P20").Value = 0.1343925781596098
$ws.Range("Q20").Value = 0.1565265643333334
$ws.Range("R20").Value = 1.408739079
$ws.Range("S20").Value = 0.000243545325968002
$ws.Range("T20").Value = 0.000243545325968002
$ws.Range("G21").Value = 0.16825
$ws.Range("H21").Value = 0.50475
$ws.Range("I21").Value = 0.001812193272151965
$ws.Range("J21").Value = 0.001812193272151965
$ws.Range("M21").Value = 0.547574
$ws.Range("N21").Value = 1.642722
$ws.Range("O21").Value = 0.07910157378579964
$ws.Range("P21").Value = 0.07910157378579964
$ws.Range("Q21").Value = 0.09212932550000001
$ws.Range("R21").Value = 0.8291639295000001
$ws.Range("S21").Value = 0.0001433473398312583
$ws.Range("T21").Value = 0.0001433473398312584
$ws.Range("G22").Value = 1.180534666666667
$ws.Range("H22").Value = 3.541604
$ws.Range("I22").Value = 0.01271534609495094
$ws.Range("J22").Value = 0.01271534609495094
$ws.Range("M22").Value = 1.780241
$ws.Range("N22").Value = 5.340723000000001
$ws.Range("O22").Value = 0.2571704734300857
$ws.Range("P22").Value = 0.2571704734300857
$ws.Range("Q22").Value = 2.101636215521334
$ws.Range("R22").Value = 18.914725939692
$ws.Range("S22").Value = 0.003270011575065925
$ws.Range("T22").Value = 0.003270011575065925
$ws.Range("G23").Value = 1.180534666666667
$ws.Range("H23").Value = 3.541604
$ws.Range("I23").Value = 0.01271534609495094
$ws.Range("J23").Value = 0.01271534609495094
$ws.Range("O23").Value = 0.4607709215973151
$ws.Range("P23").Value = 0.4607709215973152
$ws.Range("Q23").Value = 3.765490038464
$ws.Range("R23").Value = 33.889410346176
$ws.Range("S23").Value = 0.005858861738599366
$ws.Range("T23").Value = 0.005858861738599368
$ws.Range("G24").Value = 1.180534666666667
$ws.Range("H24").Value = 3.541604
$ws.Range("I24").Value = 0.01271534609495094
$ws.Range("J24").Value = 0.01271534609495094
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.4746316666666666
$ws.Range("N24").Value = 1.423895
$ws.Range("O24").Value = 0.06856445302718973
$ws.Range("P24").Value = 0.06856445302718973
$ws.Range("Q24").Value = 0.5603191363977778
$ws.Range("R24").Value = 5.042872227579999
$ws.Range("S24").Value = 0.0008718207500517241
$ws.Range("T24").Value = 0.0008718207500517242
$ws.Range("G25").Value = 1.180534666666667
$ws.Range("H25").Value = 3.541604
$ws.Range("I25").Value = 0.01271534609495094
$ws.Range("J25").Value = 0.01271534609495094
$ws.Range("M25").Value = 0.9303213333333334
$ws.Range("N25").Value = 2.790964
$ws.Range("O25").Value = 0.1343925781596098
$ws.Range("P25").Value = 0.1343925781596098
$ws.Range("Q25").Value = 1.098276585139556
$ws.Range("R25").Value = 9.884489266256001
$ws.Range("S25").Value = 0.001708848143892184
$ws.Range("T25").Value = 0.001708848143892184
$ws.Range("G26").Value = 1.180534666666667
$ws.Range("H26").Value = 3.541604
$ws.Range("I26").Value = 0.01271534609495094
$ws.Range("J26").Value = 0.01271534609495094
$ws.Range("M26").Value = 0.547574
$ws.Range("N26").Value = 1.642722
$ws.Range("O26").Value = 0.07910157378579964
$ws.Range("P26").Value = 0.07910157378579964
$ws.Range("Q26").Value = 0.6464300895653333
$ws.Range("R26").Value = 5.817870806088
$ws.Range("S26").Value = 0.001005803887341741
$ws.Range("T26").Value = 0.001005803887341741
